$wb = $excel.ActiveWorkbook

# --- Step 1: create the new "2022-Q3" sheet by copying the existing
# "2022-Q1" fund-holdings sheet (so it inherits the same header styles,
# column widths, etc.) and inserting the copy right before it, then
# overwrite its data with the new Q3 figures. ---
$wsQ1 = $wb.Worksheets.Item("2022-Q1")
$wsQ1.Copy($wsQ1, $null)
$wsQ3 = $wb.Worksheets.Item("2022-Q1 (2)")
$wsQ3.Name = "2022-Q3"

$wsQ3.Range("B2").NumberFormat = "@"
$wsQ3.Range("B2").Value = "159628"
$wsQ3.Range("C2").Value = "万家国证2000ETF"
$wsQ3.Range("D2").NumberFormat = "@"
$wsQ3.Range("D2").Value = "2.90"
$wsQ3.Range("E2").NumberFormat = "@"
$wsQ3.Range("E2").Value = "97.72"
$wsQ3.Range("F2").NumberFormat = "@"
$wsQ3.Range("F2").Value = "0.46"
$wsQ3.Range("G2").NumberFormat = "@"
$wsQ3.Range("G2").Value = "0.0133"
$wsQ3.Range("H2").Value = 7

# Drop the temporary "@" text-format now that the values are locked in as
# text, so the cells end up with no explicit style (matching the sibling
# fund sheets, which don't carry a style override on their data row).
$wsQ3.Range("B2:G2").ClearFormats()

# --- Step 2: update the "总计" summary sheet: insert a new row for
# 2022-Q3 above the existing 2022-Q1 row, pushing everything else down
# one row, and renumber the index column. ---
$wsTotal = $wb.Worksheets.Item("总计")
$wsTotal.Rows.Item(2).Insert()

# Give the new A2 the same style as the other index cells (A3/A4),
# then clear the inherited row-3 formatting off B2:D2.
$wsTotal.Range("A3").Copy()
$wsTotal.Range("A2").PasteSpecial(-4122)
$wsTotal.Range("B2:D2").ClearFormats()

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q3"
$wsTotal.Range("C2").Value = 1
$wsTotal.Range("D2").Value = 0.01

$wsTotal.Range("A3").Value = 1
$wsTotal.Range("A4").Value = 2

# Restore "总计" as the active sheet (it was the active sheet before the edit).
$wsTotal.Activate()
